$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.654.74"
$ws.Range("E2").Value = "'  -1.49%  "
$ws.Range("D3").Value = "'3.422.32"
$ws.Range("E3").Value = "'  -1.56%  "
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'572.56"
$ws.Range("E5").Value = "'  -1.24%  "
$ws.Range("D6").Value = "'157.96"
$ws.Range("E6").Value = "'  -1.56%  "
$ws.Range("D7").Value = "'0.612"
$ws.Range("E7").Value = "'  +1.39%  "
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("D9").Value = "'3.423.22"
$ws.Range("E9").Value = "'  -1.52%  "
$ws.Range("D10").Value = "'7.18"
$ws.Range("E10").Value = "'  -1.44%  "
$ws.Range("E11").Value = "'  -1.98%  "
$ws.Range("E12").Value = "'  -1.36%  "
$ws.Range("D13").Value = "'4.011.16"
$ws.Range("E13").Value = "'  -1.56%  "
$ws.Range("E14").Value = "'  -0.19%  "
$ws.Range("E15").Value = "'  -3.92%  "
$ws.Range("D16").Value = "'27.65"
$ws.Range("E16").Value = "'  -3.60%  "
$ws.Range("D17").Value = "'64.674.88"
$ws.Range("E17").Value = "'  -1.41%  "
$ws.Range("D18").Value = "'3.444.51"
$ws.Range("E18").Value = "'  -0.58%  "
$ws.Range("E19").Value = "'  -1.69%  "
$ws.Range("D20").Value = "'13.80"
$ws.Range("E20").Value = "'  -3.30%  "
$ws.Range("D21").Value = "'380.49"
$ws.Range("E21").Value = "'  -2.40%  "
$ws.Range("D22").Value = "'7.99"
$ws.Range("E22").Value = "'  -2.83%  "
$ws.Range("D23").Value = "'0.547"
$ws.Range("E23").Value = "'  -0.53%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "'  +0.12%  "
$ws.Range("D25").Value = "'72.12"
$ws.Range("E25").Value = "'  -1.75%  "
$ws.Range("E26").Value = "'  -4.61%  "
$ws.Range("D27").Value = "'9.99"
$ws.Range("E27").Value = "'  +2.34%  "
$ws.Range("D28").Value = "'0.177"
$ws.Range("E28").Value = "'  -0.24%  "
$ws.Range("D29").Value = "'1.02"
$ws.Range("E29").Value = "'  +2.00%  "
$ws.Range("E30").Value = "'  +2.75%  "
$ws.Range("E31").Value = "'  -3.60%  "
$ws.Range("D32").Value = "'1.99"
$ws.Range("E32").Value = "'  -3.09%  "
$ws.Range("D33").Value = "'23.22"
$ws.Range("E33").Value = "'  -2.03%  "
$ws.Range("D34").Value = "'7.08"
$ws.Range("E34").Value = "'  -0.13%  "
$ws.Range("D35").Value = "'1.58"
$ws.Range("E35").Value = "'  +2.41%  "
$ws.Range("D36").Value = "'160.33"
$ws.Range("E36").Value = "'  -1.77%  "
$ws.Range("E37").Value = "'  -2.42%  "
$ws.Range("B38").Value = "'Maker"
$ws.Range("C38").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "'2.906.87"
$ws.Range("E38").Value = "'  -5.52%  "
$ws.Range("B39").Value = "'Hedera"
$ws.Range("C39").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.0755"
$ws.Range("E39").Value = "'  -1.79%  "
$ws.Range("D40").Value = "'6.72"
$ws.Range("E40").Value = "'  +3.49%  "
$ws.Range("D41").Value = "'26.36"
$ws.Range("E41").Value = "'  -2.80%  "
$ws.Range("E42").Value = "'  +1.41%  "
$ws.Range("E43").Value = "'  +0.00%  "
$ws.Range("E44").Value = "'  -1.90%  "
$ws.Range("D45").Value = "'0.770"
$ws.Range("E45").Value = "'  -0.69%  "
$ws.Range("D46").Value = "'25.69"
$ws.Range("E46").Value = "'  +0.67%  "
$ws.Range("D47").Value = "'317.86"
$ws.Range("E47").Value = "'  +3.04%  "
$ws.Range("B48").Value = "'dogwifhat"
$ws.Range("C48").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.22"
$ws.Range("E48").Value = "'  -0.70%  "
$ws.Range("B49").Value = "'ONDO"
$ws.Range("C49").Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "'1.07"
$ws.Range("E49").Value = "'  -4.43%  "
$ws.Range("E50").Value = "'  -1.21%  "
$ws.Range("D51").Value = "'6.53"
$ws.Range("E51").Value = "'  -2.45%  "

# Reset style to Normal to strip the auto-applied quotePrefix formatting
# introduced by the leading apostrophe, restoring cells to their original (unstyled) state
$ws.Range("B2:E51").Style = "Normal"
